$p = $ppt.ActivePresentation

# Slide 27: add speaker notes
$s27 = $p.Slides.Item(27)
$notes27 = $s27.NotesPage.Shapes.Placeholders.Item(2)
$notes27.TextFrame.TextRange.Text = "Slightly patronising slide, but made this because I think non-coders often think having lots of error messages is a chore.  I know I used to."

# Slide 36: add speaker notes
$s36 = $p.Slides.Item(36)
$notes36 = $s36.NotesPage.Shapes.Placeholders.Item(2)
$notes36.TextFrame.TextRange.Text = "Last point will be dealt with later by Chris."
